$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, pushing the old row 22 (and anything below) down to 23
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with updated data
$ws.Range("A22").Value = 5
$ws.Range("B22").Value = "Macroferia Regional de Talca"
$ws.Range("C22").Value = "Maule"
$ws.Range("D22").Value = 44461
$ws.Range("E22").Value = 7
$ws.Range("F22").Value = 100112026
$ws.Range("G22").Value = "Haba"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 12000
$ws.Range("L22").Value = 12000
$ws.Range("M22").Value = 12000
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Provincia del Elquí"
$ws.Range("P22").Value = 480
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
